$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.215.88'
$ws.Range("E2").Value = '  -0.63%  '
$ws.Range("D3").Value = '2.269.50'
$ws.Range("E3").Value = '  -0.85%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '306.59'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.37%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.49'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.57%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.525'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.28%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  -1.27%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.32'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.40%  '
$ws.Range("E11").Value = '  -1.96%  '
$ws.Range("E12").Value = '  +0.28%  '
$ws.Range("E13").Value = '  +1.93%  '
$ws.Range("D14").Value = '2.622.41'
$ws.Range("E14").Value = '  -0.88%  '
$ws.Range("E15").Value = '  +1.20%  '
$ws.Range("D16").Value = '2.250.15'
$ws.Range("E16").Value = '  -2.48%  '
$ws.Range("E17").Value = '  -1.81%  '
$ws.Range("D18").Value = '42.120.04'
$ws.Range("E18").Value = '  -0.67%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.28'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.08%  '
$ws.Range("E20").Value = '  -1.68%  '
$ws.Range("E21").Value = '  -0.26%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.73'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.59%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '237.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.62%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.99'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.33%  '
$ws.Range("E25").Value = '  -1.20%  '
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '23.52'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.43%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '37.23'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.36%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.58'
$ws.Range("D29").Style = "Normal"
$ws.Range("E30").Value = '  +1.57%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '163.17'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.17%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.25'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.49%  '
$ws.Range("E33").Value = '  +0.04%  '
$ws.Range("E34").Value = '  +1.22%  '
$ws.Range("E35").Value = '  +2.32%  '
$ws.Range("E36").Value = '  -2.86%  '
$ws.Range("E37").Value = '  +0.16%  '
$ws.Range("E38").Value = '  -4.42%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.82'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.15%  '
$ws.Range("E40").Value = '  -1.48%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.12'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.55%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.33'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.80%  '
$ws.Range("D43").Value = '1.955.47'
$ws.Range("E43").Value = '  -2.73%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '18.99'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.08%  '
$ws.Range("E45").Value = '  -1.70%  '
$ws.Range("E46").Value = '  -1.71%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.91'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.39%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '53.74'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.29%  '
$ws.Range("D49").Value = '2.493.67'
$ws.Range("E49").Value = '  -0.54%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '72.28'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.03%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '92.54'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.08%  '
